$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The worker database rows (16-18) are reordered:
#  - Laura Vanessa Sanchez Zambrano (CC 1047452363) moves to the top,
#    with period 2407 first, then 2406.
#  - Santander Eliecer Lopez Carrillo (CC 1003367783, period 2303) moves
#    to the bottom, keeping its original values.

# Row 16: CC / 1047452363 / LAURA VANESSA SANCHEZ ZAMBRANO / 2407
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047452363"
$ws.Range("D16").Value = "LAURA VANESSA SANCHEZ ZAMBRANO"
$ws.Range("E16").Value = "2407"
$ws.Range("F16").Value = 52200
$ws.Range("G16").Value = 1305000

# Row 17: CC / 1047452363 / LAURA VANESSA SANCHEZ ZAMBRANO / 2406
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047452363"
$ws.Range("D17").Value = "LAURA VANESSA SANCHEZ ZAMBRANO"
$ws.Range("E17").Value = "2406"
$ws.Range("F17").Value = 52200
$ws.Range("G17").Value = 1305000

# Row 18: CC / 1003367783 / SANTANDER ELIECER LOPEZ CARRILLO / 2303
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1003367783"
$ws.Range("D18").Value = "SANTANDER ELIECER LOPEZ CARRILLO"
$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 21333
$ws.Range("G18").Value = 1000000
